$d = $word.ActiveDocument

# Locate the abstract paragraph (the one that currently starts with
# "For this project") rather than hard-coding an index, so the script is
# resilient to any incidental paragraph-count differences.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("For this project")) {
        $target = $p
    }
}
if ($target -eq $null) {
    Write-Output "ERROR: abstract paragraph not found"
} else {

    # --- Replace the opening of the abstract (through "...relates to ") ---
    $oldOpening = "For this project, I had to calibrate the Clio infrared detector using Python code. Due to the detector’s light saturation, there was a non-linear data trend. Therefore, I had to write code to figure out how the data should be calibrated. I first obtained data where the exposure time was gradually increased. I had to read in every picture’s exposure and brightness count, and devised a way to linearize the pictures to correct them through an equation. I took the coefficients generated and calibrated another set by applying the generated coefficients. My biggest accomplishments revolved around learning Python and seeing how coding was applied in a scientific setting. I also learned about the process of data collecting and how that relates to "
    $newOpening = "Clio is the infrared camera of the Magellan Telescope’s Adaptive Optics instrument. The Clio detector records light levels, but as those increase, the camera response isn’t linear. I measured the nonlinear response of the Clio detector and determined how to correct this. I obtained Clio images where the exposure time was gradually increased. I read in every image’s exposure and brightness count, and devised a way to correct the images by linearization through an equation. With the coefficients generated, I calibrated other data as well. My accomplishments revolved around learning Python and seeing how coding was applied in a scientific setting. I also learned about the process of data collecting and how that relates to "

    $rng1 = $target.Range
    $ok1 = $rng1.Find.Execute($oldOpening, $true, $false, $false, $false, $false, $true, 1, $false, $newOpening, 2)

    # "work beyond my undergraduate career. Since the data set is" stays untouched in between.

    # --- Replace the closing of the abstract (from " now calibrated..." to the end) ---
    $oldClosing = " now calibrated, we can fix more data sets from Clio, and use those to accurately measure the brightness of other stars and exoplanets found with this instrument."
    $newClosing = " now calibrated, we can fix more data from Clio, and use that to accurately measure the brightness of other stars and exoplanets found."

    $rng2 = $target.Range
    $ok2 = $rng2.Find.Execute($oldClosing, $true, $false, $false, $false, $false, $true, 1, $false, $newClosing, 2)

    # --- Move the _GoBack bookmark from the end of the paragraph to its start ---
    # Re-adding a bookmark with the same name over the whole (now-updated)
    # paragraph range re-anchors its start at the beginning of the range and
    # its end at the end of the range, matching the target layout.
    $d.Bookmarks.Add("_GoBack", $target.Range)

    Write-Output ("opening replaced: " + $ok1)
    Write-Output ("closing replaced: " + $ok2)
    Write-Output ("final text: " + $target.Range.Text)
}
